$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.714.51'
$ws.Range('E2').Value = '  -1.55%  '

$ws.Range('D3').Value = '3.572.30'
$ws.Range('E3').Value = '  -1.73%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').Value = '''576.49'
$ws.Range('E5').Value = '  -3.02%  '

$ws.Range('D6').Value = '''188.93'
$ws.Range('E6').Value = '  -1.99%  '

$ws.Range('D7').Value = '''0.632'
$ws.Range('E7').Value = '  -3.34%  '

$ws.Range('D8').Value = '3.568.35'
$ws.Range('E8').Value = '  -1.64%  '

$ws.Range('E9').Value = '  +0.05%  '

$ws.Range('D10').Value = '''0.178'
$ws.Range('E10').Value = '  -1.47%  '

$ws.Range('D11').Value = '''0.660'
$ws.Range('E11').Value = '  -0.78%  '

$ws.Range('D12').Value = '''55.70'
$ws.Range('E12').Value = '  -4.47%  '

$ws.Range('D13').Value = '''0.0000302'
$ws.Range('E13').Value = '  +1.89%  '

$ws.Range('D14').Value = '''9.62'
$ws.Range('E14').Value = '  -1.95%  '

$ws.Range('D15').Value = '4.143.95'
$ws.Range('E15').Value = '  -1.76%  '

$ws.Range('D16').Value = '''19.75'
$ws.Range('E16').Value = '  +1.70%  '

$ws.Range('D17').Value = '3.568.68'
$ws.Range('E17').Value = '  -1.80%  '

$ws.Range('D18').Value = '69.712.82'
$ws.Range('E18').Value = '  -1.54%  '

$ws.Range('D19').Value = '''12.63'
$ws.Range('E19').Value = '  -0.04%  '

$ws.Range('E20').Value = '  -0.14%  '

$ws.Range('D21').Value = '''1.04'
$ws.Range('E21').Value = '  -1.45%  '

$ws.Range('D22').Value = '''474.41'
$ws.Range('E22').Value = '  -4.13%  '

$ws.Range('D23').Value = '''19.46'
$ws.Range('E23').Value = '  +13.69%  '

$ws.Range('D24').Value = '''5.03'
$ws.Range('E24').Value = '  -7.62%  '

$ws.Range('D25').Value = '''4.38'
$ws.Range('E25').Value = '  -2.84%  '

$ws.Range('D26').Value = '''94.96'
$ws.Range('E26').Value = '  +4.13%  '

$ws.Range('D27').Value = '''3.01'
$ws.Range('E27').Value = '  -3.79%  '

$ws.Range('D28').Value = '''11.00'
$ws.Range('E28').Value = '  -2.58%  '

$ws.Range('D29').Value = '''9.36'
$ws.Range('E29').Value = '  -1.26%  '

$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').Value = '''7.94'
$ws.Range('E30').Value = '  +4.03%  '

$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '''32.39'
$ws.Range('E31').Value = '  -0.14%  '

$ws.Range('E32').Value = '  +1.33%  '

$ws.Range('D33').Value = '''12.21'
$ws.Range('E33').Value = '  -0.40%  '

$ws.Range('D34').Value = '''66.10'
$ws.Range('E34').Value = '  +1.13%  '

$ws.Range('D35').Value = '''580.27'
$ws.Range('E35').Value = '  -6.28%  '

$ws.Range('D36').Value = '''38.97'
$ws.Range('E36').Value = '  +1.89%  '

$ws.Range('D37').Value = '''1.00'
$ws.Range('E37').Value = '  +0.12%  '

$ws.Range('D38').Value = '0.0₃0795'
$ws.Range('E38').Value = '  -4.51%  '

$ws.Range('D39').Value = '''0.395'
$ws.Range('E39').Value = '  -4.30%  '

$ws.Range('D40').Value = '''3.20'
$ws.Range('E40').Value = '  +16.35%  '

$ws.Range('D41').Value = '''0.138'
$ws.Range('E41').Value = '  -7.28%  '

$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').Value = '''2.87'
$ws.Range('E42').Value = '  +6.17%  '

$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '''3.44'
$ws.Range('E43').Value = '  -6.18%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '3.230.11'
$ws.Range('E44').Value = '  -3.09%  '

$ws.Range('D45').Value = '''3.08'
$ws.Range('E45').Value = '  -0.12%  '

$ws.Range('D46').Value = '''0.0442'
$ws.Range('E46').Value = '  -1.59%  '

$ws.Range('D47').Value = '''3.38'
$ws.Range('E47').Value = '  +0.47%  '

$ws.Range('D48').Value = '''9.42'
$ws.Range('E48').Value = '  +1.53%  '

$ws.Range('D49').Value = '''0.138'
$ws.Range('E49').Value = '  -0.34%  '

$ws.Range('D50').Value = '''0.998'
$ws.Range('E50').Value = '  -0.22%  '

$ws.Range('D51').Value = '''3.14'
$ws.Range('E51').Value = '  -6.05%  '
